# Weekly update: insert a new week's record at the top of the data table
# (row 62) and push the existing rows 62-152 down by one (they become
# rows 63-153). This mirrors the author's commit "Fruta / hortaliza,
# semanal" where a new weekly observation is prepended to the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 62..152 down to 63..153, leaving row 62 blank (but with the
# same column formatting carried down from the row above, as Excel does).
$ws.Rows("62:62").Insert()

# Populate the newly inserted row 62 with this week's data.
$ws.Range("A62").Value = 8
$ws.Range("B62").Value = "Terminal La Palmera de La Serena"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 44763
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = 100112001
$ws.Range("G62").Value = "Berenjena"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 480
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 11000
$ws.Range("M62").Value = 10500
$ws.Range("N62").Value = "$/caja 50 unidades"
$ws.Range("O62").Value = "Región de Arica y Parinacota"
$ws.Range("P62").Value = 210
$ws.Range("Q62").Value = 50
$ws.Range("R62").Value = "Hortaliza"
